$wb = $excel.ActiveWorkbook

# Sheet names affected: "展览" and "全部类型" both need F2:F5 updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 561
    $ws.Range("F3").Value = 3568
    $ws.Range("F4").Value = 100
    $ws.Range("F5").Value = 693
}
